$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.348.12"
$ws.Range("E2").Value = "  +0.87%  "
$ws.Range("D3").Value = "3.509.51"
$ws.Range("E3").Value = "  +0.15%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "600.67"
$ws.Range("E5").Value = "  +1.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "175.58"
$ws.Range("E6").Value = "  +3.60%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  -1.04%  "
$ws.Range("E9").Value = "  -0.20%  "
$ws.Range("E10").Value = "  -2.14%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.432"
$ws.Range("E11").Value = "  -0.23%  "
$ws.Range("D12").Value = "4.114.20"
$ws.Range("E12").Value = "  +0.07%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "30.72"
$ws.Range("E13").Value = "  +8.36%  "
$ws.Range("E14").Value = "  +0.22%  "
$ws.Range("D15").Value = "67.258.36"
$ws.Range("E15").Value = "  +0.74%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000180"
$ws.Range("E16").Value = "  -1.19%  "
$ws.Range("D17").Value = "3.503.51"
$ws.Range("E17").Value = "  +0.09%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.32"
$ws.Range("E18").Value = "  -0.20%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.67"
$ws.Range("E19").Value = "  +4.46%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "395.47"
$ws.Range("E20").Value = "  -0.56%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.01"
$ws.Range("E21").Value = "  +0.34%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.41"
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.539"
$ws.Range("E24").Value = "  +0.80%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.68"
$ws.Range("E25").Value = "  -0.67%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000123"
$ws.Range("E26").Value = "  +0.30%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.20"
$ws.Range("E27").Value = "  +0.06%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.181"
$ws.Range("E28").Value = "  -0.09%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.994"
$ws.Range("E29").Value = "  -0.51%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.19"
$ws.Range("E30").Value = "  -2.32%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.43"
$ws.Range("E31").Value = "  -2.07%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.07"
$ws.Range("E32").Value = "  -0.04%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.71"
$ws.Range("E33").Value = "  -0.72%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.41"
$ws.Range("E34").Value = "  -0.06%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.64"
$ws.Range("E35").Value = "  +1.60%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "163.93"
$ws.Range("E36").Value = "  +0.81%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.879"
$ws.Range("E37").Value = "  -2.81%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.93"
$ws.Range("E38").Value = "  +0.62%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.09"
$ws.Range("E39").Value = "  +4.15%  "
$ws.Range("B40").Value = "InjectiveProtocol"
$ws.Range("C40").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "27.71"
$ws.Range("E40").Value = "  +1.59%  "
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.70"
$ws.Range("E41").Value = "  +0.16%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0734"
$ws.Range("E42").Value = "  -1.73%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "26.21"
$ws.Range("E43").Value = "  -1.35%  "
$ws.Range("D44").Value = "2.809.63"
$ws.Range("E44").Value = "  +0.07%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "42.54"
$ws.Range("E45").Value = "  -1.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.55"
$ws.Range("E46").Value = "  -0.88%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0304"
$ws.Range("E47").Value = "  -3.20%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "341.54"
$ws.Range("E48").Value = "  -0.60%  "
$ws.Range("E49").Value = "  -1.49%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "33.78"
$ws.Range("E50").Value = "  -1.08%  "
$ws.Range("B51").Value = "SuiNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.847"
$ws.Range("E51").Value = "  -1.43%  "
